# Update Betfair Back/Lay odds for the match in row 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.29   # Odd_H_Back
$ws.Range("G2").Value = 1.33   # Odd_H_Lay
$ws.Range("S2").Value = 2.56   # Odd_Over35_FT_Back
$ws.Range("W2").Value = 4      # Double_Chance_Draw_or_Away_Back
